$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is being regenerated with a re-run of the NATMI analysis
# (per Dr Hou advice): more target clusters (M1/M2 added) and refreshed
# statistics for every Sending/Target cluster combination.

# Drop the old 3 data rows first so stale strings are not left dangling
$ws.Rows("2:4").Delete()

# Sending cluster / Ligand symbol / Receptor symbol / Target cluster (text columns)
$colA = @("ECs", "ECs", "ECs", "ECs", "ECs", "FAPs", "FAPs", "FAPs", "FAPs", "FAPs")
$colB = @("Wnt5a", "Wnt5a", "Wnt5a", "Wnt5a", "Wnt5a", "Wnt5a", "Wnt5a", "Wnt5a", "Wnt5a", "Wnt5a")
$colC = @("Fzd8", "Fzd8", "Fzd8", "Fzd8", "Fzd8", "Fzd8", "Fzd8", "Fzd8", "Fzd8", "Fzd8")
$colD = @("ECs", "FAPs", "M1", "M2", "sCs", "ECs", "FAPs", "M1", "M2", "sCs")

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}

# Numeric columns E..T for rows 2-11
$numericData = @(
    @(1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 2.776364333333333, 8.329093, 0.2206324039768381, 0.2206324039768381, 0.4439415823547778, 3.995474241193, 0.005841842890697102, 0.005841842890697102),
    @(1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 8.808225333333333, 26.424676, 0.6999729490580858, 0.6999729490580857, 1.408438166875111, 12.675943501876, 0.01853368735702367, 0.01853368735702366),
    @(1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 1, 0.3333333333333333, 0.003041, 0.009122999999999999, 0.0002416624981232283, 0.0002416624981232283, 0.0004862569136666666, 0.004376312223, 0.000006398671823190071, 0.000006398671823190071),
    @(1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 1, 0.3333333333333333, 0.003928, 0.011784, 0.0003121507045800858, 0.0003121507045800858, 0.0006280885093333334, 0.005652796583999999, 0.000008265038777208352, 0.000008265038777208352),
    @(1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 0.9921066666666668, 2.97632, 0.0788408337623728, 0.07884083376237279, 0.1586381867022222, 1.42774368032, 0.002087525476356141, 0.002087525476356141),
    @(3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 2.776364333333333, 8.329093, 0.2206324039768381, 0.2206324039768381, 16.32266792304533, 146.904011307408, 0.214790561086141, 0.214790561086141),
    @(3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 8.808225333333333, 26.424676, 0.6999729490580858, 0.6999729490580857, 51.78489558491733, 466.064060264256, 0.6814392617010622, 0.6814392617010621),
    @(3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 1, 0.3333333333333333, 0.003041, 0.009122999999999999, 0.0002416624981232283, 0.0002416624981232283, 0.017878501232, 0.160906511088, 0.0002352638263000383, 0.0002352638263000383),
    @(3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 1, 0.3333333333333333, 0.003928, 0.011784, 0.0003121507045800858, 0.0003121507045800858, 0.023093309056, 0.207839781504, 0.0003038856658028775, 0.0003038856658028775),
    @(3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 0.9921066666666668, 2.97632, 0.0788408337623728, 0.07884083376237279, 5.832745893546668, 52.49471304192001, 0.07675330828601666, 0.07675330828601665),
)

for ($i = 0; $i -lt $numericData.Count; $i++) {
    $r = $i + 2
    $rowVals = $numericData[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = $j + 5
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
